$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.903.64', '  +0.00%  ', 0),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.874.11', '  -0.99%  ', 0),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  +0.03%  ', 1),
    @(5, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.7401', '  -4.12%  ', 1),
    @(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '242.30', '  -0.72%  ', 1),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9998', '  -0.08%  ', 1),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3146', '  +0.66%  ', 1),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07169', '  -0.87%  ', 1),
    @(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.70', '  -3.84%  ', 1),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.08442', '  -2.90%  ', 1),
    @(12, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7522', '  -2.45%  ', 1),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.397', '  -0.07%  ', 1),
    @(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.871.78', '  -7.16%  ', 0),
    @(15, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '92.55', '  -1.78%  ', 1),
    @(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.903.15', '  -0.72%  ', 0),
    @(17, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.095', '  -1.69%  ', 1),
    @(18, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.58', '  -2.46%  ', 1),
    @(19, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '243.59', '  -0.67%  ', 1),
    @(20, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007817', '  -0.51%  ', 1),
    @(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9992', '  -0.19%  ', 1),
    @(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.126.64', '  -8.51%  ', 0),
    @(23, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.991', '  -2.12%  ', 1),
    @(24, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  +0.04%  ', 1),
    @(25, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1556', '  -2.25%  ', 1),
    @(26, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.320', '  -2.11%  ', 1),
    @(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '165.63', '  +2.00%  ', 1),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.61', '  -1.06%  ', 1),
    @(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.042', '  +0.06%  ', 1),
    @(30, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.478', '  +3.03%  ', 1),
    @(31, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.610', '  +2.08%  ', 1),
    @(32, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.530', '  -0.91%  ', 1),
    @(33, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.256', '  +3.44%  ', 1),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05326', '  -2.58%  ', 1),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.241', '  -0.62%  ', 1),
    @(36, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7559', '  +0.37%  ', 1),
    @(37, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9974', '  -0.98%  ', 1),
    @(38, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.696', '  -0.18%  ', 1),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01955', '  -0.77%  ', 1),
    @(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.748', '  -1.37%  ', 1),
    @(41, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4490', '  -0.41%  ', 1),
    @(42, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.111.92', '  +1.42%  ', 0),
    @(43, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.081', '  +0.00%  ', 1),
    @(44, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '72.31', '  -1.74%  ', 1),
    @(45, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8586', '  +0.49%  ', 1),
    @(46, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.000', '  +0.01%  ', 1),
    @(47, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '103.14', '  +0.09%  ', 1),
    @(48, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.673', '  +0.72%  ', 1),
    @(49, 'SynthetixNetwork', 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx', '3.068', '  +2.92%  ', 1),
    @(50, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.840', '  -2.21%  ', 1),
    @(51, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.023.57', '  -8.76%  ', 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[5] -eq 1) {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
